$wb = $excel.ActiveWorkbook

# --- Add "4Y Blocks Data" sheet after "11Y Blocks" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "4Y Blocks Data"
$ws3.PageSetup.LeftMargin = 54
$ws3.PageSetup.RightMargin = 54
$ws3.PageSetup.TopMargin = 72
$ws3.PageSetup.BottomMargin = 72
$ws3.PageSetup.HeaderMargin = 36
$ws3.PageSetup.FooterMargin = 36

# Header row
$ws3.Cells.Item(1,1).Value = "1980-1983"
$ws3.Cells.Item(1,2).Value = "1984-1987"
$ws3.Cells.Item(1,3).Value = "1988-1991"
$ws3.Cells.Item(1,4).Value = "1992-1995"
$ws3.Cells.Item(1,5).Value = "1996-1999"
$ws3.Cells.Item(1,6).Value = "2000-2003"
$ws3.Cells.Item(1,7).Value = "2004-2007"
$ws3.Cells.Item(1,8).Value = "2008-2011"
$ws3.Cells.Item(1,9).Value = "2012-2015"
$ws3.Cells.Item(1,10).Value = "2016-2019"
$ws3.Cells.Item(1,11).Value = "2020-2023"

# Data rows
$ws3.Cells.Item(2,1).Value = 1984215.04
$ws3.Cells.Item(2,2).Value = 1885275.517
$ws3.Cells.Item(2,3).Value = 1923381.114
$ws3.Cells.Item(2,4).Value = 1777557.862
$ws3.Cells.Item(2,5).Value = 1862215.666
$ws3.Cells.Item(2,6).Value = 1683637.982
$ws3.Cells.Item(2,7).Value = 1900290.624
$ws3.Cells.Item(2,8).Value = 1790393.904
$ws3.Cells.Item(2,9).Value = 1815850.831
$ws3.Cells.Item(2,10).Value = 1601410.129
$ws3.Cells.Item(2,11).Value = 1754165.795
$ws3.Cells.Item(3,1).Value = 1838153.82
$ws3.Cells.Item(3,2).Value = 1940873.569
$ws3.Cells.Item(3,3).Value = 1799819.844
$ws3.Cells.Item(3,4).Value = 1970400.155
$ws3.Cells.Item(3,5).Value = 1973296.457
$ws3.Cells.Item(3,6).Value = 1783958.153
$ws3.Cells.Item(3,7).Value = 1826529.309
$ws3.Cells.Item(3,8).Value = 1826130.974
$ws3.Cells.Item(3,9).Value = 1641965.98
$ws3.Cells.Item(3,10).Value = 1700190.457
$ws3.Cells.Item(3,11).Value = 1646240.327
$ws3.Cells.Item(4,1).Value = 2068803.876
$ws3.Cells.Item(4,2).Value = 1951143.692
$ws3.Cells.Item(4,3).Value = 1942424.567
$ws3.Cells.Item(4,4).Value = 1990792.359
$ws3.Cells.Item(4,5).Value = 1633907.377
$ws3.Cells.Item(4,6).Value = 1925298.682
$ws3.Cells.Item(4,7).Value = 1842572.17
$ws3.Cells.Item(4,8).Value = 1441530.109
$ws3.Cells.Item(4,9).Value = 1781967.197
$ws3.Cells.Item(4,10).Value = 1813670.267
$ws3.Cells.Item(4,11).Value = 1513700.956
$ws3.Cells.Item(5,1).Value = 1926352.241
$ws3.Cells.Item(5,2).Value = 2105383.706
$ws3.Cells.Item(5,3).Value = 2001966.894
$ws3.Cells.Item(5,4).Value = 1790723.418
$ws3.Cells.Item(5,5).Value = 1713535.308
$ws3.Cells.Item(5,6).Value = 1883313.057
$ws3.Cells.Item(5,7).Value = 1786635.778
$ws3.Cells.Item(5,8).Value = 1763011.252
$ws3.Cells.Item(5,9).Value = 1895202.71
$ws3.Cells.Item(5,10).Value = 1949969.043
$ws3.Cells.Item(5,11).Value = 1846196.197

# --- Add "11Y Blocks Data" sheet after "4Y Blocks Data" ---
$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "11Y Blocks Data"
$ws4.PageSetup.LeftMargin = 54
$ws4.PageSetup.RightMargin = 54
$ws4.PageSetup.TopMargin = 72
$ws4.PageSetup.BottomMargin = 72
$ws4.PageSetup.HeaderMargin = 36
$ws4.PageSetup.FooterMargin = 36

# Header row
$ws4.Cells.Item(1,1).Value = "1980-1990"
$ws4.Cells.Item(1,2).Value = "1991-2001"
$ws4.Cells.Item(1,3).Value = "2002-2012"
$ws4.Cells.Item(1,4).Value = "2013-2023"

# Data rows
$ws4.Cells.Item(2,1).Value = 1984215.04
$ws4.Cells.Item(2,2).Value = 2001966.894
$ws4.Cells.Item(2,3).Value = 1925298.682
$ws4.Cells.Item(2,4).Value = 1641965.98
$ws4.Cells.Item(3,1).Value = 1838153.82
$ws4.Cells.Item(3,2).Value = 1777557.862
$ws4.Cells.Item(3,3).Value = 1883313.057
$ws4.Cells.Item(3,4).Value = 1781967.197
$ws4.Cells.Item(4,1).Value = 2068803.876
$ws4.Cells.Item(4,2).Value = 1970400.155
$ws4.Cells.Item(4,3).Value = 1900290.624
$ws4.Cells.Item(4,4).Value = 1895202.71
$ws4.Cells.Item(5,1).Value = 1926352.241
$ws4.Cells.Item(5,2).Value = 1990792.359
$ws4.Cells.Item(5,3).Value = 1826529.309
$ws4.Cells.Item(5,4).Value = 1601410.129
$ws4.Cells.Item(6,1).Value = 1885275.517
$ws4.Cells.Item(6,2).Value = 1790723.418
$ws4.Cells.Item(6,3).Value = 1842572.17
$ws4.Cells.Item(6,4).Value = 1700190.457
$ws4.Cells.Item(7,1).Value = 1940873.569
$ws4.Cells.Item(7,2).Value = 1862215.666
$ws4.Cells.Item(7,3).Value = 1786635.778
$ws4.Cells.Item(7,4).Value = 1813670.267
$ws4.Cells.Item(8,1).Value = 1951143.692
$ws4.Cells.Item(8,2).Value = 1973296.457
$ws4.Cells.Item(8,3).Value = 1790393.904
$ws4.Cells.Item(8,4).Value = 1949969.043
$ws4.Cells.Item(9,1).Value = 2105383.706
$ws4.Cells.Item(9,2).Value = 1633907.377
$ws4.Cells.Item(9,3).Value = 1826130.974
$ws4.Cells.Item(9,4).Value = 1754165.795
$ws4.Cells.Item(10,1).Value = 1923381.114
$ws4.Cells.Item(10,2).Value = 1713535.308
$ws4.Cells.Item(10,3).Value = 1441530.109
$ws4.Cells.Item(10,4).Value = 1646240.327
$ws4.Cells.Item(11,1).Value = 1799819.844
$ws4.Cells.Item(11,2).Value = 1683637.982
$ws4.Cells.Item(11,3).Value = 1763011.252
$ws4.Cells.Item(11,4).Value = 1513700.956
$ws4.Cells.Item(12,1).Value = 1942424.567
$ws4.Cells.Item(12,2).Value = 1783958.153
$ws4.Cells.Item(12,3).Value = 1815850.831
$ws4.Cells.Item(12,4).Value = 1846196.197

# --- Apply header formatting (bold, centered, bordered) matching existing sheets ---
$ws1 = $wb.Worksheets.Item("4Y Blocks")
$ws1.Range("A1").Copy()
$ws3.Range("A1:K1").PasteSpecial(-4122)
$ws1.Range("A1").Copy()
$ws4.Range("A1:D1").PasteSpecial(-4122)

